$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.644.93"
$ws.Range("E2").Value = "'  +2.25%  "
$ws.Range("D3").Value = "'1.892.70"
$ws.Range("E4").Value = "'  +0.16%  "
$ws.Range("D5").Value = "'244.38"
$ws.Range("E5").Value = "'  +1.35%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "'  +0.12%  "
$ws.Range("D7").Value = "'0.4964"
$ws.Range("E7").Value = "'  +0.10%  "
$ws.Range("D8").Value = "'0.2963"
$ws.Range("E8").Value = "'  +1.99%  "
$ws.Range("D9").Value = "'0.06814"
$ws.Range("E9").Value = "'  +3.61%  "
$ws.Range("D10").Value = "'1.892.79"
$ws.Range("E10").Value = "'  +1.03%  "
$ws.Range("D11").Value = "'17.10"
$ws.Range("E11").Value = "'  +2.43%  "
$ws.Range("E12").Value = "'  +2.11%  "
$ws.Range("D13").Value = "'91.19"
$ws.Range("E13").Value = "'  +6.11%  "
$ws.Range("D14").Value = "'5.099"
$ws.Range("E14").Value = "'  +5.41%  "
$ws.Range("D15").Value = "'0.6748"
$ws.Range("E15").Value = "'  +2.50%  "
$ws.Range("D16").Value = "'30.645.04"
$ws.Range("E16").Value = "'  +2.30%  "
$ws.Range("D17").Value = "'0.000007934"
$ws.Range("E17").Value = "'  +0.87%  "
$ws.Range("E18").Value = "'  +0.19%  "
$ws.Range("D19").Value = "'13.28"
$ws.Range("E19").Value = "'  +4.78%  "
$ws.Range("D20").Value = "'2.136.32"
$ws.Range("E20").Value = "'  +1.07%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "'  +0.11%  "
$ws.Range("D22").Value = "'4.864"
$ws.Range("E22").Value = "'  +2.83%  "
$ws.Range("D23").Value = "'179.60"
$ws.Range("E23").Value = "'  +33.27%  "
$ws.Range("D24").Value = "'6.058"
$ws.Range("E24").Value = "'  +8.59%  "
$ws.Range("D25").Value = "'9.302"
$ws.Range("E25").Value = "'  +3.03%  "
$ws.Range("D26").Value = "'154.02"
$ws.Range("E26").Value = "'  +2.58%  "
$ws.Range("D27").Value = "'18.73"
$ws.Range("E27").Value = "'  +12.29%  "
$ws.Range("D28").Value = "'1.929"
$ws.Range("E28").Value = "'  +1.68%  "
$ws.Range("D29").Value = "'1.386"
$ws.Range("E29").Value = "'  +1.09%  "
$ws.Range("D30").Value = "'4.334"
$ws.Range("E30").Value = "'  +4.51%  "
$ws.Range("D31").Value = "'0.08930"
$ws.Range("E31").Value = "'  +2.90%  "
$ws.Range("D32").Value = "'4.040"
$ws.Range("E32").Value = "'  +2.82%  "
$ws.Range("D33").Value = "'0.05206"
$ws.Range("E33").Value = "'  +3.56%  "
$ws.Range("D34").Value = "'0.7374"
$ws.Range("E34").Value = "'  +5.55%  "
$ws.Range("E35").Value = "'  +3.81%  "
$ws.Range("D36").Value = "'2.672"
$ws.Range("E36").Value = "'  +0.51%  "
$ws.Range("D37").Value = "'0.01878"
$ws.Range("E37").Value = "'  +10.65%  "
$ws.Range("D38").Value = "'2.701"
$ws.Range("E38").Value = "'  +0.59%  "
$ws.Range("D39").Value = "'2.170"
$ws.Range("E39").Value = "'  +0.41%  "
$ws.Range("D40").Value = "'0.9344"
$ws.Range("E40").Value = "'  +1.02%  "
$ws.Range("D41").Value = "'0.4362"
$ws.Range("E41").Value = "'  +4.50%  "
$ws.Range("D42").Value = "'106.01"
$ws.Range("E42").Value = "'  +4.37%  "
$ws.Range("D43").Value = "'5.802"
$ws.Range("E43").Value = "'  -2.20%  "
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "'  +0.27%  "
$ws.Range("D45").Value = "'7.656"
$ws.Range("E45").Value = "'  +3.94%  "
$ws.Range("D47").Value = "'0.05845"
$ws.Range("E47").Value = "'  +3.63%  "
$ws.Range("D48").Value = "'33.38"
$ws.Range("D49").Value = "'0.3896"
$ws.Range("D50").Value = "'8.501"
$ws.Range("E50").Value = "'  +5.20%  "
$ws.Range("E51").Value = "'  +3.79%  "
